$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sdvfvvfdvd"
$ws.Range("A2").Value = "vzfvd"
$ws.Range("A3").Value = "vdf"

$ws.Range("A3").Select()
